# Update the "dSF" column (column F) with repulled data.
# Column E ("dS0") is unchanged; column F ("dSF") gets new values
# reflecting the freshly pulled data / recalculated means.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newF = @{
    2  = -2
    3  = 6
    4  = 2
    5  = 3
    6  = -6
    7  = -1
    8  = 1
    9  = -1
    10 = 1
    11 = -1
    12 = -1
    13 = 0
    14 = 5
    15 = 1
    16 = 0
    17 = -1
    18 = -1
    19 = 4
    20 = 0
    21 = 0
    22 = -2
    23 = 2
    24 = 5
    25 = -4
    26 = -2
    27 = -1
    28 = 5
    29 = 2
    30 = 4
    31 = 3
    32 = 2
    33 = -6
    34 = 1
    35 = 0
    36 = 3
    37 = -2
    38 = 3
    39 = 0
    40 = 0
    41 = 0
}

foreach ($row in $newF.Keys) {
    $ws.Range("F$row").Value = $newF[$row]
}
